$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Tiempo de Ejecución" column (G) entirely
$ws.Columns("G").Delete()

# Update existing rows and add new iter3 rows (A:F)
$ws.Cells.Item(2, 1).Value = "new_values_log"
$ws.Cells.Item(2, 2).Value = "iter3-UnderSampling"
$ws.Cells.Item(2, 3).Value = "FL"
$ws.Cells.Item(2, 4).Value = "RandomForest"
$ws.Cells.Item(2, 5).Value = "{'model__n_estimators': 50}"
$ws.Cells.Item(2, 6).Value = 0.8469430984397408

$ws.Cells.Item(3, 1).Value = "new_values_log"
$ws.Cells.Item(3, 2).Value = "iter3-OverSampling"
$ws.Cells.Item(3, 3).Value = "FL"
$ws.Cells.Item(3, 4).Value = "RandomForest"
$ws.Cells.Item(3, 5).Value = "{'model__n_estimators': 50}"
$ws.Cells.Item(3, 6).Value = 0.8376728270197041

$ws.Cells.Item(4, 1).Value = "new_values_log"
$ws.Cells.Item(4, 2).Value = "iter3-UnderSampling"
$ws.Cells.Item(4, 3).Value = "FL"
$ws.Cells.Item(4, 4).Value = "RandomForest"
$ws.Cells.Item(4, 5).Value = "{'model__n_estimators': 50}"
$ws.Cells.Item(4, 6).Value = 0.8358860457707065

$ws.Cells.Item(5, 1).Value = "new_values"
$ws.Cells.Item(5, 2).Value = "iter3-OverSampling"
$ws.Cells.Item(5, 3).Value = "GA"
$ws.Cells.Item(5, 4).Value = "Lasso"
$ws.Cells.Item(5, 5).Value = "{'model__alpha': 1.0}"
$ws.Cells.Item(5, 6).Value = 0.9001475754583788

$ws.Cells.Item(6, 1).Value = "new_values"
$ws.Cells.Item(6, 2).Value = "iter3-UnderSampling"
$ws.Cells.Item(6, 3).Value = "GA"
$ws.Cells.Item(6, 4).Value = "Lasso"
$ws.Cells.Item(6, 5).Value = "{'model__alpha': 1.0}"
$ws.Cells.Item(6, 6).Value = 0.9001258546622485

$ws.Cells.Item(7, 1).Value = "new_values"
$ws.Cells.Item(7, 2).Value = "iter2-cov"
$ws.Cells.Item(7, 3).Value = "GA"
$ws.Cells.Item(7, 4).Value = "Lasso"
$ws.Cells.Item(7, 5).Value = "{'model__alpha': 1.0}"
$ws.Cells.Item(7, 6).Value = 0.8788511235177607

$ws.Cells.Item(8, 1).Value = "new_values"
$ws.Cells.Item(8, 2).Value = "iter3-OverSampling"
$ws.Cells.Item(8, 3).Value = "NC"
$ws.Cells.Item(8, 4).Value = "LinearRegression"
$ws.Cells.Item(8, 5).Value = "{}"
$ws.Cells.Item(8, 6).Value = 0.9176460433793464

$ws.Cells.Item(9, 1).Value = "new_values"
$ws.Cells.Item(9, 2).Value = "iter3-UnderSampling"
$ws.Cells.Item(9, 3).Value = "NC"
$ws.Cells.Item(9, 4).Value = "LinearRegression"
$ws.Cells.Item(9, 5).Value = "{}"
$ws.Cells.Item(9, 6).Value = 0.9176460433793464

$ws.Cells.Item(10, 1).Value = "new_values"
$ws.Cells.Item(10, 2).Value = "iter2-cov"
$ws.Cells.Item(10, 3).Value = "NC"
$ws.Cells.Item(10, 4).Value = "LinearRegression"
$ws.Cells.Item(10, 5).Value = "{}"
$ws.Cells.Item(10, 6).Value = 0.9119573949198786

$ws.Cells.Item(11, 1).Value = "all_log"
$ws.Cells.Item(11, 2).Value = "iter3-OverSampling"
$ws.Cells.Item(11, 3).Value = "NJ"
$ws.Cells.Item(11, 4).Value = "AdaBoost"
$ws.Cells.Item(11, 5).Value = "{'model__n_estimators': 50}"
$ws.Cells.Item(11, 6).Value = 0.9017978960530233

$ws.Cells.Item(12, 1).Value = "all_log"
$ws.Cells.Item(12, 2).Value = "iter3-UnderSampling"
$ws.Cells.Item(12, 3).Value = "NJ"
$ws.Cells.Item(12, 4).Value = "AdaBoost"
$ws.Cells.Item(12, 5).Value = "{'model__n_estimators': 50}"
$ws.Cells.Item(12, 6).Value = 0.7997092373207666

$ws.Cells.Item(13, 1).Value = "all_log"
$ws.Cells.Item(13, 2).Value = "iter2-cov"
$ws.Cells.Item(13, 3).Value = "NJ"
$ws.Cells.Item(13, 4).Value = "AdaBoost"
$ws.Cells.Item(13, 5).Value = "{'model__n_estimators': 50}"
$ws.Cells.Item(13, 6).Value = 0.6143037383604859

$ws.Cells.Item(14, 1).Value = "new_values"
$ws.Cells.Item(14, 2).Value = "iter3-UnderSampling"
$ws.Cells.Item(14, 3).Value = "NY"
$ws.Cells.Item(14, 4).Value = "GradientBoosting"
$ws.Cells.Item(14, 5).Value = "{'model__n_estimators': 150}"
$ws.Cells.Item(14, 6).Value = 0.9529815385623106

$ws.Cells.Item(15, 1).Value = "new_values"
$ws.Cells.Item(15, 2).Value = "iter3-OverSampling"
$ws.Cells.Item(15, 3).Value = "NY"
$ws.Cells.Item(15, 4).Value = "GradientBoosting"
$ws.Cells.Item(15, 5).Value = "{'model__n_estimators': 150}"
$ws.Cells.Item(15, 6).Value = 0.9403092362887392

$ws.Cells.Item(16, 1).Value = "new_values"
$ws.Cells.Item(16, 2).Value = "iter1-normal"
$ws.Cells.Item(16, 3).Value = "NY"
$ws.Cells.Item(16, 4).Value = "GradientBoosting"
$ws.Cells.Item(16, 5).Value = "{'model__n_estimators': 150}"
$ws.Cells.Item(16, 6).Value = 0.9060651223899857

$ws.Cells.Item(17, 1).Value = "new_values"
$ws.Cells.Item(17, 2).Value = "iter3-OverSampling"
$ws.Cells.Item(17, 3).Value = "SC"
$ws.Cells.Item(17, 4).Value = "Lasso"
$ws.Cells.Item(17, 5).Value = "{'model__alpha': 1.0}"
$ws.Cells.Item(17, 6).Value = 0.8458742226569919

$ws.Cells.Item(18, 1).Value = "new_values"
$ws.Cells.Item(18, 2).Value = "iter3-UnderSampling"
$ws.Cells.Item(18, 3).Value = "SC"
$ws.Cells.Item(18, 4).Value = "Lasso"
$ws.Cells.Item(18, 5).Value = "{'model__alpha': 1.0}"
$ws.Cells.Item(18, 6).Value = 0.8458742226569919

$ws.Cells.Item(19, 1).Value = "new_values"
$ws.Cells.Item(19, 2).Value = "iter2-cov"
$ws.Cells.Item(19, 3).Value = "SC"
$ws.Cells.Item(19, 4).Value = "Lasso"
$ws.Cells.Item(19, 5).Value = "{'model__alpha': 1.0}"
$ws.Cells.Item(19, 6).Value = 0.8355560181287938

$ws.Cells.Item(20, 1).Value = "new_values"
$ws.Cells.Item(20, 2).Value = "iter3-OverSampling"
$ws.Cells.Item(20, 3).Value = "VA"
$ws.Cells.Item(20, 4).Value = "LinearRegression"
$ws.Cells.Item(20, 5).Value = "{}"
$ws.Cells.Item(20, 6).Value = 0.9409347961886814

$ws.Cells.Item(21, 1).Value = "new_values"
$ws.Cells.Item(21, 2).Value = "iter3-UnderSampling"
$ws.Cells.Item(21, 3).Value = "VA"
$ws.Cells.Item(21, 4).Value = "LinearRegression"
$ws.Cells.Item(21, 5).Value = "{}"
$ws.Cells.Item(21, 6).Value = 0.9409347961886814

$ws.Cells.Item(22, 1).Value = "new_values"
$ws.Cells.Item(22, 2).Value = "iter2-cov"
$ws.Cells.Item(22, 3).Value = "VA"
$ws.Cells.Item(22, 4).Value = "LinearRegression"
$ws.Cells.Item(22, 5).Value = "{}"
$ws.Cells.Item(22, 6).Value = 0.921118736676871
